$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure price/volume columns that look numeric stay stored as text,
# matching the inlineStr cell type used in the source workbook.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "69.302.57"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "3.898.53"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "528.66"
$ws.Range("D6").Value = "144.65"
$ws.Range("E6").Value = "  -0.60%  "
$ws.Range("E7").Value = "  -1.80%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("E10").Value = "  -4.50%  "
$ws.Range("E11").Value = "  -5.29%  "
$ws.Range("D12").Value = "42.09"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "4.517.64"
$ws.Range("E13").Value = "  +0.25%  "
$ws.Range("D14").Value = "10.25"
$ws.Range("E14").Value = "  -2.53%  "
$ws.Range("D15").Value = "3.920.01"
$ws.Range("E15").Value = "  +1.00%  "
$ws.Range("E16").Value = "  -2.01%  "
$ws.Range("E17").Value = "  -1.26%  "
$ws.Range("E18").Value = "  +6.68%  "
$ws.Range("E19").Value = "  -0.83%  "
$ws.Range("D20").Value = "69.267.38"
$ws.Range("E20").Value = "  +1.51%  "
$ws.Range("D21").Value = "426.02"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("E22").Value = "  -5.54%  "
$ws.Range("D23").Value = "14.14"
$ws.Range("E23").Value = "  -4.38%  "
$ws.Range("D24").Value = "88.10"
$ws.Range("E24").Value = "  -0.83%  "
$ws.Range("D25").Value = "4.04"
$ws.Range("E25").Value = "  +10.27%  "
$ws.Range("D26").Value = "11.39"
$ws.Range("E26").Value = "  -10.07%  "
$ws.Range("D27").Value = "10.61"
$ws.Range("E27").Value = "  -3.69%  "
$ws.Range("D28").Value = "36.37"
$ws.Range("E28").Value = "  -2.33%  "
$ws.Range("D29").Value = "688.67"
$ws.Range("E29").Value = "  -4.38%  "
$ws.Range("D30").Value = "13.18"
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("E31").Value = "  -3.12%  "
$ws.Range("D32").Value = "2.81"
$ws.Range("E32").Value = "  -3.61%  "
$ws.Range("D33").Value = "68.88"
$ws.Range("E33").Value = "  +11.53%  "
$ws.Range("D34").Value = "0.0₃0884"
$ws.Range("E34").Value = "  +0.56%  "
$ws.Range("D35").Value = "0.436"
$ws.Range("E35").Value = "  +8.92%  "
$ws.Range("E36").Value = "  -2.03%  "
$ws.Range("D37").Value = "39.97"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("E38").Value = "  +2.14%  "
$ws.Range("E39").Value = "  +0.06%  "
$ws.Range("E40").Value = "  +7.53%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  +9.06%  "
$ws.Range("E43").Value = "  -3.48%  "
$ws.Range("E44").Value = "  -6.35%  "
$ws.Range("D45").Value = "3.41"
$ws.Range("E45").Value = "  +1.13%  "
$ws.Range("E46").Value = "  +18.79%  "
$ws.Range("E47").Value = "  -1.55%  "
$ws.Range("D48").Value = "2.99"
$ws.Range("E48").Value = "  +6.83%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "147.95"
$ws.Range("E49").Value = "  +2.34%  "
$ws.Range("B50").Value = "BabyDogeCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D50").Value = "0.0₆0346"
$ws.Range("E50").Value = "  -3.78%  "
$ws.Range("D51").Value = "2.743.70"
$ws.Range("E51").Value = "  +14.69%  "

Write-Output "Applied cryptos update"
